$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 22:44"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 4231344
$ws.Range("C4").Value = 61026
$ws.Range("D4").Value = 1999212
$ws.Range("E4").Value = 2083948
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 835
$ws.Range("H4").Value = 148184

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1337022
$ws.Range("C6").Value = 48892
$ws.Range("D6").Value = 850107
$ws.Range("E6").Value = 455509
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 761
$ws.Range("H6").Value = 31406

# Row 11: Chile
$ws.Range("A11").Value = "Chile"
$ws.Range("B11").Value = 341304
$ws.Range("C11").Value = 2545
$ws.Range("D11").Value = 313696
$ws.Range("E11").Value = 18694
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 76
$ws.Range("H11").Value = 8914

# Row 21: Alemania
$ws.Range("A21").Value = "Alemania"
$ws.Range("B21").Value = 205879
$ws.Range("C21").Value = 737
$ws.Range("D21").Value = 189400
$ws.Range("E21").Value = 7286
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 9193

# Row 28: Egipto
$ws.Range("A28").Value = "Egipto"
$ws.Range("B28").Value = 91072
$ws.Range("C28").Value = 659
$ws.Range("D28").Value = 31970
$ws.Range("E28").Value = 54584
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 38
$ws.Range("H28").Value = 4518

# Row 30: Ecuador
$ws.Range("A30").Value = "Ecuador"
$ws.Range("B30").Value = 79049
$ws.Range("C30").Value = 901
$ws.Range("D30").Value = 34544
$ws.Range("E30").Value = 39037
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 29
$ws.Range("H30").Value = 5468

# Row 31: Suecia
$ws.Range("A31").Value = "Suecia"
$ws.Range("B31").Value = 78763
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 5676

# Row 32: Kazajistan
$ws.Range("A32").Value = "Kazajistan"
$ws.Range("B32").Value = 78486
$ws.Range("C32").Value = 1687
$ws.Range("D32").Value = 49488
$ws.Range("E32").Value = 28413
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 585

# Row 93: Gabon
$ws.Range("A93").Value = "Gabon"
$ws.Range("B93").Value = 6984
$ws.Range("C93").Value = 396
$ws.Range("D93").Value = 4463
$ws.Range("E93").Value = 2472
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 49

# Row 94: Guinea
$ws.Range("A94").Value = "Guinea"
$ws.Range("B94").Value = 6867
$ws.Range("C94").Value = 61
$ws.Range("D94").Value = 6063
$ws.Range("E94").Value = 762
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 42

# Row 99: Republica de Africa Central
$ws.Range("A99").Value = "Republica de Africa Central"
$ws.Range("B99").Value = 4593
$ws.Range("C99").Value = 3
$ws.Range("D99").Value = 1483
$ws.Range("E99").Value = 3051
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 59

# Row 131: Ruanda
$ws.Range("A131").Value = "Ruanda"
$ws.Range("B131").Value = 1729
$ws.Range("C131").Value = 19
$ws.Range("D131").Value = 900
$ws.Range("E131").Value = 824
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 5

# Row 172: Martinica
$ws.Range("A172").Value = "Martinica"
$ws.Range("B172").Value = 269
$ws.Range("C172").Value = 7
$ws.Range("D172").Value = 98
$ws.Range("E172").Value = 156
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 15

# Row 174: Gambia
$ws.Range("A174").Value = "Gambia"
$ws.Range("B174").Value = 216
$ws.Range("C174").Value = 46
$ws.Range("D174").Value = 60
$ws.Range("E174").Value = 150
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 6

# Row 175: Islas Caimanes
$ws.Range("A175").Value = "Islas Caimanes"
$ws.Range("B175").Value = 203
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 202
$ws.Range("E175").Value = 0
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 1

# Row 176: Camboya
$ws.Range("A176").Value = "Camboya"
$ws.Range("B176").Value = 202
$ws.Range("C176").Value = 4
$ws.Range("D176").Value = 142
$ws.Range("E176").Value = 60
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

# Row 177: Guadalupe
$ws.Range("A177").Value = "Guadalupe"
$ws.Range("B177").Value = 195
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 172
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 14

# Row 178: Islas Feroe
$ws.Range("A178").Value = "Islas Feroe"
$ws.Range("B178").Value = 191
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 188
$ws.Range("E178").Value = 3
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179: Gibraltar
$ws.Range("A179").Value = "Gibraltar"
$ws.Range("B179").Value = 184
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 180
$ws.Range("E179").Value = 4
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

# Row 210: Groenlandia
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("B210").Value = 13
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 13
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Islas Malvinas
$ws.Range("A211").Value = "Islas Malvinas"
$ws.Range("B211").Value = 13
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 13
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

Write-Output "done"